$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 4. Источники данных и методы сбора -> Источники данных (B16):
# replace the old "ф.№3-здрав" sourcing text with the updated "ф.№14" text.
$ws.Range("B16").Value = 'Источниками информации для расчета показателя является административные данные Минздрава КР. На ежегодной основе собирается государственная форма отчетности ф.№14 "Отчет о деятельности стационара".'

# 4. Источники данных и методы сбора -> Методы сбора данных (B17):
# replace the old "№3-ЗДРАВ" collection-method text with the updated "№14" text.
$ws.Range("B17").Value = 'Областные медико-информационные центры осуществляют сбор статистической отчетности по форме №14 "Отчет о деятельности стационара". Метод наблюдения – сплошные наблюдения. Единицами наблюдения являются каждые роды.'

# Leave the selection on the last-edited cell, matching the author's final view.
$ws.Range("B16").Select()
